# Apply updated crypto price/volume data per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '57.535.16'
$ws.Range('E2').Value2 = '  -4.52%  '
$ws.Range('D3').Value2 = '3.082.15'
$ws.Range('E3').Value2 = '  -6.54%  '
$ws.Range('E4').Value2 = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '514.20'
$ws.Range('E5').Value2 = '  -7.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '128.80'
$ws.Range('E6').Value2 = '  -8.57%  '
$ws.Range('E7').Value2 = '  +0.18%  '
$ws.Range('D8').Value2 = '3.069.50'
$ws.Range('E8').Value2 = '  -6.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.435'
$ws.Range('E9').Value2 = '  -6.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '7.13'
$ws.Range('E10').Value2 = '  -9.60%  '
$ws.Range('E11').Value2 = '  -11.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.364'
$ws.Range('E12').Value2 = '  -10.56%  '
$ws.Range('D13').Value2 = '3.621.76'
$ws.Range('E13').Value2 = '  -6.20%  '
$ws.Range('E14').Value2 = '  -0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '24.45'
$ws.Range('E15').Value2 = '  -7.76%  '
$ws.Range('D16').Value2 = '57.708.57'
$ws.Range('E16').Value2 = '  -4.20%  '
$ws.Range('D17').Value2 = '3.093.63'
$ws.Range('E17').Value2 = '  -6.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '0.0000147'
$ws.Range('E18').Value2 = '  -10.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '5.57'
$ws.Range('E19').Value2 = '  -7.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '12.65'
$ws.Range('E20').Value2 = '  -7.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '7.64'
$ws.Range('E21').Value2 = '  -10.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '333.29'
$ws.Range('E22').Value2 = '  -10.72%  '
$ws.Range('E23').Value2 = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '0.498'
$ws.Range('E24').Value2 = '  -6.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '65.81'
$ws.Range('E25').Value2 = '  -8.52%  '
$ws.Range('E26').Value2 = '  -4.20%  '
$ws.Range('E27').Value2 = '  +0.20%  '
$ws.Range('D28').Value2 = '0.0₃0891'
$ws.Range('E28').Value2 = '  -12.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '1.00'
$ws.Range('E29').Value2 = '  +0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '6.62'
$ws.Range('E30').Value2 = '  -6.14%  '
$ws.Range('E32').Value2 = '  -10.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '6.66'
$ws.Range('E33').Value2 = '  -8.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '20.80'
$ws.Range('E34').Value2 = '  -7.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '157.48'
$ws.Range('E35').Value2 = '  -4.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '4.66'
$ws.Range('E36').Value2 = '  -7.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '5.98'
$ws.Range('E37').Value2 = '  -9.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '1.33'
$ws.Range('E38').Value2 = '  -12.50%  '
$ws.Range('D39').Value2 = '3.128.70'
$ws.Range('E39').Value2 = '  -6.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '39.94'
$ws.Range('E40').Value2 = '  -4.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.0664'
$ws.Range('E41').Value2 = '  -7.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '22.72'
$ws.Range('E42').Value2 = '  -10.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.671'
$ws.Range('E43').Value2 = '  -10.02%  '
$ws.Range('E44').Value2 = '  +0.46%  '
$ws.Range('B45').Value2 = 'ONDO'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '1.04'
$ws.Range('E45').Value2 = '  -6.83%  '
$ws.Range('B46').Value2 = 'Filecoin'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '3.80'
$ws.Range('E46').Value2 = '  -7.25%  '
$ws.Range('D47').Value2 = '2.235.51'
$ws.Range('E47').Value2 = '  -3.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '1.39'
$ws.Range('E48').Value2 = '  -11.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '6.00'
$ws.Range('E49').Value2 = '  -5.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '19.92'
$ws.Range('E50').Value2 = '  -7.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '0.0227'
$ws.Range('E51').Value2 = '  -9.83%  '
